$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper scratch cell (off the used range) for writing percentage strings
# like "84%" as literal text: assigning such a string straight to a target
# cell via .Value is auto-parsed by Excel into the number 0.84 with a new
# percent number format. Building it as a text formula in the scratch cell
# and pasting only the resulting value preserves both the literal text and
# the destination cell's original style.

$ws.Range("E2").Value = "2026-02-09 19:18:27"
$ws.Range("O2").Value = "-0.8 °C"
$ws.Range("E3").Value = "2026-02-09 19:18:30"
$ws.Range("E4").Value = "2026-02-09 19:18:32"
$ws.Range("ZZ1").Formula = "=""84%"""
$ws.Range("ZZ1").Copy()
$ws.Range("H4").PasteSpecial(-4163)
$ws.Range("E5").Value = "2026-02-09 19:18:34"
$ws.Range("O5").Value = "-2.7 °C"
$ws.Range("E6").Value = "2026-02-09 19:18:37"
$ws.Range("ZZ1").Formula = "=""84%"""
$ws.Range("ZZ1").Copy()
$ws.Range("H6").PasteSpecial(-4163)
$ws.Range("E7").Value = "2026-02-09 19:18:39"
$ws.Range("E8").Value = "2026-02-09 19:18:42"
$ws.Range("ZZ1").Formula = "=""69%"""
$ws.Range("ZZ1").Copy()
$ws.Range("H8").PasteSpecial(-4163)
$ws.Range("E9").Value = "2026-02-09 19:18:44"
$ws.Range("ZZ1").Formula = "=""80%"""
$ws.Range("ZZ1").Copy()
$ws.Range("H9").PasteSpecial(-4163)
$ws.Range("E10").Value = "2026-02-09 19:18:47"
$ws.Range("E11").Value = "2026-02-09 19:18:49"
$ws.Range("O11").Value = "5.4 °C"
$ws.Range("E12").Value = "2026-02-09 19:18:51"
$ws.Range("ZZ1").Formula = "=""85%"""
$ws.Range("ZZ1").Copy()
$ws.Range("H12").PasteSpecial(-4163)
$ws.Range("E13").Value = "2026-02-09 19:18:53"
$ws.Range("O13").Value = "3.0 °C"
$ws.Range("E14").Value = "2026-02-09 19:18:56"
$ws.Range("E15").Value = "2026-02-09 19:18:58"
$ws.Range("E16").Value = "2026-02-09 19:19:01"
$ws.Range("ZZ1").Formula = "=""72%"""
$ws.Range("ZZ1").Copy()
$ws.Range("H16").PasteSpecial(-4163)
$ws.Range("O16").Value = "-3.5 °C"
$ws.Range("E17").Value = "2026-02-09 19:19:03"
$ws.Range("O17").Value = "1.0 °C"
$ws.Range("E18").Value = "2026-02-09 19:19:05"
$ws.Range("ZZ1").Formula = "=""81%"""
$ws.Range("ZZ1").Copy()
$ws.Range("H18").PasteSpecial(-4163)
$ws.Range("E19").Value = "2026-02-09 19:19:08"
$ws.Range("E20").Value = "2026-02-09 19:19:10"
$ws.Range("E21").Value = "2026-02-09 19:19:13"
$ws.Range("J21").Value = "1007.8 hPa"
$ws.Range("O21").Value = "4.6 °C"
$ws.Range("E22").Value = "2026-02-09 19:19:15"
$ws.Range("O22").Value = "-5.1 °C"
$ws.Range("E23").Value = "2026-02-09 19:19:17"
$ws.Range("E24").Value = "2026-02-09 19:19:20"
$ws.Range("E25").Value = "2026-02-09 19:19:22"
$ws.Range("E26").Value = "2026-02-09 19:19:24"
$ws.Range("E27").Value = "2026-02-09 19:19:27"
$ws.Range("E28").Value = "2026-02-09 19:19:29"
$ws.Range("E29").Value = "2026-02-09 19:19:32"
$ws.Range("E30").Value = "2026-02-09 19:19:34"
$ws.Range("ZZ1").Formula = "=""85%"""
$ws.Range("ZZ1").Copy()
$ws.Range("H30").PasteSpecial(-4163)
$ws.Range("J30").Value = "1007.0 hPa"
$ws.Range("E31").Value = "2026-02-09 19:19:37"
$ws.Range("J31").Value = "1006.4 hPa"
$ws.Range("E32").Value = "2026-02-09 19:19:39"
$ws.Range("O32").Value = "5.1 °C"
$ws.Range("E33").Value = "2026-02-09 19:19:41"
$ws.Range("E34").Value = "2026-02-09 19:19:44"
$ws.Range("E35").Value = "2026-02-09 19:19:46"
$ws.Range("O35").Value = "5.4 °C"
$ws.Range("E36").Value = "2026-02-09 19:19:49"
$ws.Range("ZZ1").Formula = "=""78%"""
$ws.Range("ZZ1").Copy()
$ws.Range("H36").PasteSpecial(-4163)
$ws.Range("O36").Value = "9.8 °C"
$ws.Range("E37").Value = "2026-02-09 19:19:51"
$ws.Range("O37").Value = "5.5 °C"
$ws.Range("E38").Value = "2026-02-09 19:19:53"
$ws.Range("E39").Value = "2026-02-09 19:19:56"
$ws.Range("E40").Value = "2026-02-09 19:19:58"
$ws.Range("J40").Value = "1008.6 hPa"
$ws.Range("O40").Value = "4.7 °C"
$ws.Range("E41").Value = "2026-02-09 19:20:00"
$ws.Range("E42").Value = "2026-02-09 19:20:03"
$ws.Range("E43").Value = "2026-02-09 19:20:05"
$ws.Range("E44").Value = "2026-02-09 19:20:07"
$ws.Range("E45").Value = "2026-02-09 19:20:10"
$ws.Range("ZZ1").Formula = "=""84%"""
$ws.Range("ZZ1").Copy()
$ws.Range("H45").PasteSpecial(-4163)
$ws.Range("O45").Value = "3.9 °C"
$ws.Range("E46").Value = "2026-02-09 19:20:12"
$ws.Range("ZZ1").Formula = "=""72%"""
$ws.Range("ZZ1").Copy()
$ws.Range("H46").PasteSpecial(-4163)
$ws.Range("I46").Value = "0.5 mm"
$ws.Range("O46").Value = "10.4 °C"

# Clean up the scratch cell
$ws.Range("ZZ1").Clear()
